$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Bmpr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.668521
$ws.Range("H2").Value = 5.005563
$ws.Range("I2").Value = 0.9677024783929865
$ws.Range("J2").Value = 0.9677024783929865
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.918859
$ws.Range("N2").Value = 29.756577
$ws.Range("O2").Value = 0.160764128269069
$ws.Range("P2").Value = 0.160764128269069
$ws.Range("Q2").Value = 16.549824537539
$ws.Range("R2").Value = 148.948420837851
$ws.Range("S2").Value = 0.1555718453626661
$ws.Range("T2").Value = 0.1555718453626661

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Bmpr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.668521
$ws.Range("H3").Value = 5.005563
$ws.Range("I3").Value = 0.9677024783929865
$ws.Range("J3").Value = 0.9677024783929865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 29.20351433333333
$ws.Range("N3").Value = 87.610543
$ws.Range("O3").Value = 0.47332838627826
$ws.Range("P3").Value = 0.4733283862782601
$ws.Range("Q3").Value = 48.72667693896767
$ws.Range("R3").Value = 438.5400924507091
$ws.Range("S3").Value = 0.4580410524952251
$ws.Range("T3").Value = 0.4580410524952251

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Bmpr1a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.668521
$ws.Range("H4").Value = 5.005563
$ws.Range("I4").Value = 0.9677024783929865
$ws.Range("J4").Value = 0.9677024783929865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.575837
$ws.Range("N4").Value = 67.72751099999999
$ws.Range("O4").Value = 0.3659074854526709
$ws.Range("P4").Value = 0.3659074854526709
$ws.Range("Q4").Value = 37.668258127077
$ws.Range("R4").Value = 339.014323143693
$ws.Range("S4").Value = 0.3540895805350953
$ws.Range("T4").Value = 0.3540895805350953

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Bmpr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05568766666666666
$ws.Range("H5").Value = 0.167063
$ws.Range("I5").Value = 0.03229752160701353
$ws.Range("J5").Value = 0.03229752160701353
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.918859
$ws.Range("N5").Value = 29.756577
$ws.Range("O5").Value = 0.160764128269069
$ws.Range("P5").Value = 0.160764128269069
$ws.Range("Q5").Value = 0.5523581137056666
$ws.Range("R5").Value = 4.971223023351
$ws.Range("S5").Value = 0.005192282906402952
$ws.Range("T5").Value = 0.005192282906402952

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Bmpr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05568766666666666
$ws.Range("H6").Value = 0.167063
$ws.Range("I6").Value = 0.03229752160701353
$ws.Range("J6").Value = 0.03229752160701353
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 29.20351433333333
$ws.Range("N6").Value = 87.610543
$ws.Range("O6").Value = 0.47332838627826
$ws.Range("P6").Value = 0.4733283862782601
$ws.Range("Q6").Value = 1.626275571689889
$ws.Range("R6").Value = 14.636480145209
$ws.Range("S6").Value = 0.01528733378303495
$ws.Range("T6").Value = 0.01528733378303495

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Bmpr1a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05568766666666666
$ws.Range("H7").Value = 0.167063
$ws.Range("I7").Value = 0.03229752160701353
$ws.Range("J7").Value = 0.03229752160701353
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.575837
$ws.Range("N7").Value = 67.72751099999999
$ws.Range("O7").Value = 0.3659074854526709
$ws.Range("P7").Value = 0.3659074854526709
$ws.Range("Q7").Value = 1.257195685577
$ws.Range("R7").Value = 11.314761170193
$ws.Range("S7").Value = 0.01181790491757563
$ws.Range("T7").Value = 0.01181790491757563

